$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diccionario_validador")

# Add new validation data in column E for rows 12, 13 and 14
$ws.Range("E12").Value = "proceso.val_num_11"
$ws.Range("E13").Value = "proceso.val_en_dominio"
$ws.Range("E14").Value = "proceso.val_en_dominio"

# Match style/format of other cells in column E (same style as E11, s="2")
$ws.Range("E12").NumberFormat = $ws.Range("E11").NumberFormat
$ws.Range("E13").NumberFormat = $ws.Range("E11").NumberFormat
$ws.Range("E14").NumberFormat = $ws.Range("E11").NumberFormat

$ws.Range("E12").Font.Name = $ws.Range("E11").Font.Name
$ws.Range("E13").Font.Name = $ws.Range("E11").Font.Name
$ws.Range("E14").Font.Name = $ws.Range("E11").Font.Name

$ws.Range("E12").Font.Size = $ws.Range("E11").Font.Size
$ws.Range("E13").Font.Size = $ws.Range("E11").Font.Size
$ws.Range("E14").Font.Size = $ws.Range("E11").Font.Size

# Update the selected cell to match the recorded user selection
$ws.Range("E22").Select()
